$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 453, shifting existing rows 453..545 down to 454..546
$ws.Rows.Item(453).Insert()

# Populate the newly inserted row 453 with the new data record
$ws.Cells.Item(453, 1).Value = 1
$ws.Cells.Item(453, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(453, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(453, 4).Value = 45244
$ws.Cells.Item(453, 5).Value = 15
$ws.Cells.Item(453, 6).Value = "Fruta"
$ws.Cells.Item(453, 7).Value = 100102
$ws.Cells.Item(453, 8).Value = "Cítricos"
$ws.Cells.Item(453, 9).Value = 100102003
$ws.Cells.Item(453, 10).Value = "Limón"
$ws.Cells.Item(453, 11).Value = "Sin especificar"
$ws.Cells.Item(453, 12).Value = "2a amarillo"
$ws.Cells.Item(453, 13).Value = 300
$ws.Cells.Item(453, 14).Value = 20000
$ws.Cells.Item(453, 15).Value = 22000
$ws.Cells.Item(453, 16).Value = 21000
$ws.Cells.Item(453, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(453, 18).Value = "Región Metropolitana"
$ws.Cells.Item(453, 19).Value = 1050
$ws.Cells.Item(453, 20).Value = 20

# Ensure the date cell (column D) carries the same date style as the rest of the column
$ws.Cells.Item(453, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
